$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (26 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 641.1818
$ws.Range("I99").Value = 619.25
$ws.Range("J99").Value = 699.6667
$ws.Range("K99").Value = 1857.75
$ws.Range("L99").Value = 2099.0001
$ws.Range("M99").Value = -359.75
$ws.Range("N99").Value = -5095.0001
$ws.Range("H132").Value = 5525042
$ws.Range("I132").Value = 6244917
$ws.Range("K132").Value = 18734751
$ws.Range("M132").Value = -18732221
$ws.Range("H137").Value = 30517.166
$ws.Range("I137").Value = 47957.855
$ws.Range("J137").Value = 6100.2
$ws.Range("K137").Value = 143873.565
$ws.Range("L137").Value = 18300.6
$ws.Range("M137").Value = -141323.565
$ws.Range("N137").Value = -23400.6
$ws.Range("H138").Value = 1386.3024
$ws.Range("I138").Value = 1082.1052
$ws.Range("K138").Value = 3246.3156
$ws.Range("M138").Value = 1893.6844
$ws.Range("H141").Value = 2076.9
$ws.Range("I141").Value = 1509.4
$ws.Range("K141").Value = 4528.200000000001
$ws.Range("M141").Value = 651.7999999999993

# --- Sheet: ARM (44 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16577.754
$ws.Range("I32").Value = 17560.854
$ws.Range("J32").Value = 3207.6
$ws.Range("K32").Value = 17560.854
$ws.Range("L32").Value = 3207.6
$ws.Range("M32").Value = -17273.854
$ws.Range("N32").Value = -3781.6
$ws.Range("H45").Value = 5400.8
$ws.Range("I45").Value = 2501.75
$ws.Range("K45").Value = 2501.75
$ws.Range("M45").Value = -2124.75
$ws.Range("H74").Value = 278747.38
$ws.Range("I74").Value = 376109.44
$ws.Range("J74").Value = 19115.166
$ws.Range("K74").Value = 376109.44
$ws.Range("L74").Value = 19115.166
$ws.Range("M74").Value = -375235.44
$ws.Range("N74").Value = -20863.166
$ws.Range("H77").Value = 278747.38
$ws.Range("I77").Value = 376109.44
$ws.Range("J77").Value = 19115.166
$ws.Range("K77").Value = 1880547.2
$ws.Range("L77").Value = 95575.83
$ws.Range("M77").Value = -1876179.2
$ws.Range("N77").Value = -104311.83
$ws.Range("H109").Value = 176665
$ws.Range("J109").Value = 176665
$ws.Range("L109").Value = 176665
$ws.Range("N109").Value = -179439
$ws.Range("H122").Value = 2135
$ws.Range("I122").Value = 2017.0769
$ws.Range("J122").Value = 2305.3333
$ws.Range("K122").Value = 6051.2307
$ws.Range("L122").Value = 6915.999899999999
$ws.Range("M122").Value = -3601.2307
$ws.Range("N122").Value = -11815.9999
$ws.Range("H132").Value = 2452.8572
$ws.Range("I132").Value = 1778.4166
$ws.Range("K132").Value = 5335.2498
$ws.Range("M132").Value = -2805.2498
$ws.Range("H135").Value = 99999.5
$ws.Range("J135").Value = 99999.5
$ws.Range("L135").Value = 99999.5
$ws.Range("N135").Value = -110139.5

# --- Sheet: BSM (23 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 32777.777
$ws.Range("J63").Value = 32777.777
$ws.Range("L63").Value = 32777.777
$ws.Range("N63").Value = -34149.777
$ws.Range("H66").Value = 32777.777
$ws.Range("J66").Value = 32777.777
$ws.Range("L66").Value = 98333.33100000001
$ws.Range("N66").Value = -105197.331
$ws.Range("H107").Value = 28348.895
$ws.Range("J107").Value = 4654.2
$ws.Range("L107").Value = 4654.2
$ws.Range("N107").Value = -8494.200000000001
$ws.Range("H134").Value = 3930.3264
$ws.Range("I134").Value = 2407.25
$ws.Range("J134").Value = 8148.077
$ws.Range("K134").Value = 7221.75
$ws.Range("L134").Value = 24444.231
$ws.Range("M134").Value = -4686.75
$ws.Range("N134").Value = -29514.231
$ws.Range("H135").Value = 82365.57000000001
$ws.Range("J135").Value = 81093.164
$ws.Range("L135").Value = 81093.164
$ws.Range("N135").Value = -91233.164

# --- Sheet: CRP (29 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2175162.8
$ws.Range("I31").Value = 2382107
$ws.Range("J31").Value = 2248
$ws.Range("K31").Value = 2382107
$ws.Range("L31").Value = 2248
$ws.Range("M31").Value = -2381812
$ws.Range("N31").Value = -2838
$ws.Range("H34").Value = 2175162.8
$ws.Range("I34").Value = 2382107
$ws.Range("J34").Value = 2248
$ws.Range("K34").Value = 2382107
$ws.Range("L34").Value = 2248
$ws.Range("M34").Value = -2381905
$ws.Range("N34").Value = -2652
$ws.Range("H122").Value = 2146.0833
$ws.Range("I122").Value = 1994.1
$ws.Range("K122").Value = 5982.299999999999
$ws.Range("M122").Value = -3532.299999999999
$ws.Range("H134").Value = 1967.7567
$ws.Range("I134").Value = 1794.5143
$ws.Range("K134").Value = 5383.5429
$ws.Range("M134").Value = -2848.5429
$ws.Range("H140").Value = 110797.18
$ws.Range("I140").Value = 59000
$ws.Range("J140").Value = 115976.9
$ws.Range("K140").Value = 59000
$ws.Range("L140").Value = 115976.9
$ws.Range("M140").Value = -53820
$ws.Range("N140").Value = -126336.9

# --- Sheet: CUL (16 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2983.3333
$ws.Range("I5").Value = 3475
$ws.Range("K5").Value = 10425
$ws.Range("M5").Value = -10313
$ws.Range("H23").Value = 381.6
$ws.Range("J23").Value = 616
$ws.Range("L23").Value = 1848
$ws.Range("N23").Value = -2318
$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665
$ws.Range("H135").Value = 2983.3333
$ws.Range("I135").Value = 3475
$ws.Range("K135").Value = 31275
$ws.Range("M135").Value = -28740

# --- Sheet: GSM (48 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 16484.666
$ws.Range("J24").Value = 16695.654
$ws.Range("L24").Value = 16695.654
$ws.Range("N24").Value = -17041.654
$ws.Range("H49").Value = 40002.332
$ws.Range("J49").Value = 40002.332
$ws.Range("L49").Value = 40002.332
$ws.Range("N49").Value = -40370.332
$ws.Range("H70").Value = 4249.75
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 5999.5
$ws.Range("K70").Value = 2500
$ws.Range("L70").Value = 5999.5
$ws.Range("M70").Value = -2230
$ws.Range("N70").Value = -6539.5
$ws.Range("H73").Value = 4249.75
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 5999.5
$ws.Range("K73").Value = 2500
$ws.Range("L73").Value = 5999.5
$ws.Range("M73").Value = -1564
$ws.Range("N73").Value = -7871.5
$ws.Range("H102").Value = 19180.234
$ws.Range("I102").Value = 24542.924
$ws.Range("K102").Value = 24542.924
$ws.Range("M102").Value = -22920.924
$ws.Range("H122").Value = 3525.7058
$ws.Range("I122").Value = 3494.1667
$ws.Range("K122").Value = 10482.5001
$ws.Range("M122").Value = -8032.500100000001
$ws.Range("H126").Value = 4092.2
$ws.Range("I126").Value = 2081.8
$ws.Range("J126").Value = 6102.6
$ws.Range("K126").Value = 6245.400000000001
$ws.Range("L126").Value = 18307.8
$ws.Range("M126").Value = -3775.400000000001
$ws.Range("N126").Value = -23247.8
$ws.Range("H132").Value = 4966.1665
$ws.Range("I132").Value = 5759.4
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 17278.2
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -14748.2
$ws.Range("N132").Value = -8060
$ws.Range("H136").Value = 11142.883
$ws.Range("J136").Value = 11142.883
$ws.Range("L136").Value = 33428.649
$ws.Range("N136").Value = -38528.649

# --- Sheet: LTW (31 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 33069
$ws.Range("J42").Value = 33069
$ws.Range("L42").Value = 33069
$ws.Range("N42").Value = -34195
$ws.Range("H49").Value = 33069
$ws.Range("J49").Value = 33069
$ws.Range("L49").Value = 33069
$ws.Range("N49").Value = -33363
$ws.Range("H93").Value = 2880.32
$ws.Range("I93").Value = 2045.2941
$ws.Range("K93").Value = 2045.2941
$ws.Range("M93").Value = -797.2941000000001
$ws.Range("H100").Value = 6528.067
$ws.Range("J100").Value = 15798.333
$ws.Range("L100").Value = 15798.333
$ws.Range("N100").Value = -16880.333
$ws.Range("H127").Value = 142666
$ws.Range("J127").Value = 142666
$ws.Range("L127").Value = 142666
$ws.Range("N127").Value = -152586
$ws.Range("H132").Value = 2019.6976
$ws.Range("I132").Value = 2001.7632
$ws.Range("K132").Value = 6005.2896
$ws.Range("M132").Value = -3475.2896
$ws.Range("H136").Value = 2834.1738
$ws.Range("I136").Value = 2632.4285
$ws.Range("J136").Value = 4952.5
$ws.Range("K136").Value = 7897.2855
$ws.Range("L136").Value = 14857.5
$ws.Range("M136").Value = -5347.2855
$ws.Range("N136").Value = -19957.5

# --- Sheet: WVR (14 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18439.436
$ws.Range("I132").Value = 21268.87
$ws.Range("K132").Value = 63806.61
$ws.Range("M132").Value = -61276.61
$ws.Range("H136").Value = 20856.354
$ws.Range("I136").Value = 27905.584
$ws.Range("J136").Value = 3938.2
$ws.Range("K136").Value = 83716.75199999999
$ws.Range("L136").Value = 11814.6
$ws.Range("M136").Value = -81166.75199999999
$ws.Range("H137").Value = 86744.75
$ws.Range("J137").Value = 93791.60000000001
$ws.Range("L137").Value = 93791.60000000001
$ws.Range("N137").Value = -103991.6

Write-Host "Applied 231 cell updates across 8 sheets"